$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value = "'3.81"
$ws.Range("B12").Value = "'44.48"
$ws.Range("C12").Value = "'52.99"
$ws.Range("D12").Value = "'97.46"
